$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.475.03"
$ws.Range("D3").Value = "2.372.42"
$ws.Range("E3").Value = "  +4.78%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.657"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.58"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +13.96%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "27.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").Value = "2.730.62"
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "2.368.07"
$ws.Range("E18").Value = "  +4.53%  "
$ws.Range("D19").Value = "43.471.29"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").Value = "0.0₃0995"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  +15.70%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.47%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("E32").Value = "  -7.99%  "
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  +8.20%  "
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "1.443.67"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "2.599.16"
$ws.Range("E50").Value = "  +5.03%  "
$ws.Range("E51").Value = "  -0.10%  "
